$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("SFIA Level") to make room for
# the new "Skill Description" column. This shifts the old SFIA Level,
# Keycode and Description columns from B,C,D to C,D,E.
$ws.Columns.Item(2).Insert()

# New header for inserted column B
$ws.Range("B1").Value = "Skill Description"

# Fill in the "Skill Description" value for every data row (2-14).
# For most skill groups the description equals the skill code itself;
# for PROF it expands to the full descriptive phrase.
$descriptions = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "PROF"       = "Portfolio, programme and project support"
    "MADE"       = "MADE"
}

for ($row = 2; $row -le 14; $row++) {
    $code = $ws.Range("A$row").Text
    $ws.Range("B$row").Value = $descriptions[$code]
}
